# Auto-generated edit script applying the diff to all 3 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:52:26"
$ws.Cells.Item(3,1).Value = "Total filas: 116"
$ws.Cells.Item(44,1).Value = "06:52:31"
$ws.Cells.Item(44,2).Value = "07:05"
$ws.Cells.Item(44,3).Value = "15_ABASTO"
$ws.Cells.Item(44,4).Value = 13
$ws.Cells.Item(44,5).Value = "LP1912"
$ws.Cells.Item(45,1).Value = "06:52:31"
$ws.Cells.Item(45,2).Value = "07:05"
$ws.Cells.Item(45,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(45,4).Value = 13
$ws.Cells.Item(45,5).Value = "LP1912"
$ws.Cells.Item(56,1).Value = "07:17:57"
$ws.Cells.Item(56,2).Value = "07:31"
$ws.Cells.Item(56,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(56,4).Value = 14
$ws.Cells.Item(56,5).Value = "LP1912"
$ws.Cells.Item(57,1).Value = "07:17:57"
$ws.Cells.Item(57,2).Value = "07:31"
$ws.Cells.Item(57,3).Value = "16_SANTA ANA"
$ws.Cells.Item(57,4).Value = 14
$ws.Cells.Item(57,5).Value = "LP1912"
$ws.Cells.Item(75,1).Value = "07:50:23"
$ws.Cells.Item(75,2).Value = "08:22"
$ws.Cells.Item(75,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(75,4).Value = 32
$ws.Cells.Item(75,5).Value = "LP1912"
$ws.Cells.Item(76,1).Value = "07:17:57"
$ws.Cells.Item(76,2).Value = "08:22"
$ws.Cells.Item(76,3).Value = "215B_EL PATO"
$ws.Cells.Item(76,4).Value = 65
$ws.Cells.Item(76,5).Value = "LP1912"
$ws.Cells.Item(85,1).Value = "07:50:23"
$ws.Cells.Item(85,2).Value = "08:53"
$ws.Cells.Item(85,3).Value = "10_OLMOS"
$ws.Cells.Item(85,4).Value = 63
$ws.Cells.Item(85,5).Value = "LP1912"
$ws.Cells.Item(86,1).Value = "07:17:57"
$ws.Cells.Item(86,2).Value = "08:53"
$ws.Cells.Item(86,3).Value = "17_ROMERO"
$ws.Cells.Item(86,4).Value = 96
$ws.Cells.Item(86,5).Value = "LP1912"
$ws.Cells.Item(87,1).Value = "08:52:26"
$ws.Cells.Item(87,2).Value = "08:54"
$ws.Cells.Item(87,3).Value = "17_ROMERO"
$ws.Cells.Item(87,4).Value = 2
$ws.Cells.Item(87,5).Value = "LP1912"
$ws.Cells.Item(88,1).Value = "08:52:26"
$ws.Cells.Item(88,2).Value = "08:54"
$ws.Cells.Item(88,3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(88,4).Value = 2
$ws.Cells.Item(88,5).Value = "LP1912"
$ws.Cells.Item(89,1).Value = "08:52:26"
$ws.Cells.Item(89,2).Value = "08:55"
$ws.Cells.Item(89,3).Value = "10_OLMOS"
$ws.Cells.Item(89,4).Value = 3
$ws.Cells.Item(89,5).Value = "LP1912"
$ws.Cells.Item(90,1).Value = "08:52:26"
$ws.Cells.Item(90,2).Value = "09:01"
$ws.Cells.Item(90,3).Value = "215A_EL PATO"
$ws.Cells.Item(90,4).Value = 9
$ws.Cells.Item(90,5).Value = "LP1912"
$ws.Cells.Item(91,1).Value = "08:52:26"
$ws.Cells.Item(91,2).Value = "09:03"
$ws.Cells.Item(91,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(91,4).Value = 11
$ws.Cells.Item(91,5).Value = "LP1912"
$ws.Cells.Item(92,1).Value = "08:52:26"
$ws.Cells.Item(92,2).Value = "09:04"
$ws.Cells.Item(92,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(92,4).Value = 12
$ws.Cells.Item(92,5).Value = "LP1912"
$ws.Cells.Item(93,1).Value = "08:16:28"
$ws.Cells.Item(93,2).Value = "09:08"
$ws.Cells.Item(93,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(93,4).Value = 52
$ws.Cells.Item(93,5).Value = "LP1912"
$ws.Cells.Item(94,1).Value = "08:52:26"
$ws.Cells.Item(94,2).Value = "09:10"
$ws.Cells.Item(94,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(94,4).Value = 18
$ws.Cells.Item(94,5).Value = "LP1912"
$ws.Cells.Item(95,1).Value = "08:16:28"
$ws.Cells.Item(95,2).Value = "09:13"
$ws.Cells.Item(95,3).Value = "10_OLMOS"
$ws.Cells.Item(95,4).Value = 57
$ws.Cells.Item(95,5).Value = "LP1912"
$ws.Cells.Item(96,1).Value = "08:52:26"
$ws.Cells.Item(96,2).Value = "09:16"
$ws.Cells.Item(96,3).Value = "27_EL RETIRO"
$ws.Cells.Item(96,4).Value = 24
$ws.Cells.Item(96,5).Value = "LP1912"
$ws.Cells.Item(97,1).Value = "07:50:23"
$ws.Cells.Item(97,2).Value = "09:17"
$ws.Cells.Item(97,3).Value = "27_EL RETIRO"
$ws.Cells.Item(97,4).Value = 87
$ws.Cells.Item(97,5).Value = "LP1912"
$ws.Cells.Item(98,1).Value = "08:52:26"
$ws.Cells.Item(98,2).Value = "09:21"
$ws.Cells.Item(98,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(98,4).Value = 29
$ws.Cells.Item(98,5).Value = "LP1912"
$ws.Cells.Item(99,1).Value = "08:52:26"
$ws.Cells.Item(99,2).Value = "09:21"
$ws.Cells.Item(99,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99,4).Value = 29
$ws.Cells.Item(99,5).Value = "LP1912"
$ws.Cells.Item(100,1).Value = "08:52:26"
$ws.Cells.Item(100,2).Value = "09:22"
$ws.Cells.Item(100,3).Value = "16_SANTA ANA"
$ws.Cells.Item(100,4).Value = 30
$ws.Cells.Item(100,5).Value = "LP1912"
$ws.Cells.Item(101,1).Value = "08:39:38"
$ws.Cells.Item(101,2).Value = "09:22"
$ws.Cells.Item(101,3).Value = "17_ROMERO"
$ws.Cells.Item(101,4).Value = 43
$ws.Cells.Item(101,5).Value = "LP1912"
$ws.Cells.Item(102,1).Value = "07:50:23"
$ws.Cells.Item(102,2).Value = "09:23"
$ws.Cells.Item(102,3).Value = "17_ROMERO"
$ws.Cells.Item(102,4).Value = 93
$ws.Cells.Item(102,5).Value = "LP1912"
$ws.Cells.Item(103,1).Value = "08:52:26"
$ws.Cells.Item(103,2).Value = "09:23"
$ws.Cells.Item(103,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(103,4).Value = 31
$ws.Cells.Item(103,5).Value = "LP1912"
$ws.Cells.Item(104,1).Value = "08:16:28"
$ws.Cells.Item(104,2).Value = "09:29"
$ws.Cells.Item(104,3).Value = "16_SANTA ANA"
$ws.Cells.Item(104,4).Value = 73
$ws.Cells.Item(104,5).Value = "LP1912"
$ws.Cells.Item(105,1).Value = "07:50:23"
$ws.Cells.Item(105,2).Value = "09:31"
$ws.Cells.Item(105,3).Value = "16_SANTA ANA"
$ws.Cells.Item(105,4).Value = 101
$ws.Cells.Item(105,5).Value = "LP1912"
$ws.Cells.Item(106,1).Value = "08:52:26"
$ws.Cells.Item(106,2).Value = "09:32"
$ws.Cells.Item(106,3).Value = "15_ABASTO"
$ws.Cells.Item(106,4).Value = 40
$ws.Cells.Item(106,5).Value = "LP1912"
$ws.Cells.Item(107,1).Value = "08:52:26"
$ws.Cells.Item(107,2).Value = "09:33"
$ws.Cells.Item(107,3).Value = "10_OLMOS"
$ws.Cells.Item(107,4).Value = 41
$ws.Cells.Item(107,5).Value = "LP1912"
$ws.Cells.Item(108,1).Value = "08:52:26"
$ws.Cells.Item(108,2).Value = "09:34"
$ws.Cells.Item(108,3).Value = "16_SANTA ANA"
$ws.Cells.Item(108,4).Value = 42
$ws.Cells.Item(108,5).Value = "LP1912"
$ws.Cells.Item(109,1).Value = "08:39:38"
$ws.Cells.Item(109,2).Value = "09:39"
$ws.Cells.Item(109,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(109,4).Value = 60
$ws.Cells.Item(109,5).Value = "LP1912"
$ws.Cells.Item(110,1).Value = "08:39:38"
$ws.Cells.Item(110,2).Value = "09:41"
$ws.Cells.Item(110,3).Value = "215C_EL PATO"
$ws.Cells.Item(110,4).Value = 62
$ws.Cells.Item(110,5).Value = "LP1912"
$ws.Cells.Item(111,1).Value = "08:39:38"
$ws.Cells.Item(111,2).Value = "09:42"
$ws.Cells.Item(111,3).Value = "10_OLMOS"
$ws.Cells.Item(111,4).Value = 63
$ws.Cells.Item(111,5).Value = "LP1912"
$ws.Cells.Item(112,1).Value = "08:52:26"
$ws.Cells.Item(112,2).Value = "09:42"
$ws.Cells.Item(112,3).Value = "215C_EL PATO"
$ws.Cells.Item(112,4).Value = 50
$ws.Cells.Item(112,5).Value = "LP1912"
$ws.Cells.Item(113,1).Value = "08:52:26"
$ws.Cells.Item(113,2).Value = "09:43"
$ws.Cells.Item(113,3).Value = "14_ABASTO"
$ws.Cells.Item(113,4).Value = 51
$ws.Cells.Item(113,5).Value = "LP1912"
$ws.Cells.Item(114,1).Value = "08:52:26"
$ws.Cells.Item(114,2).Value = "09:52"
$ws.Cells.Item(114,3).Value = "15_ABASTO"
$ws.Cells.Item(114,4).Value = 60
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Cells.Item(115,1).Value = "08:52:26"
$ws.Cells.Item(115,2).Value = "09:53"
$ws.Cells.Item(115,3).Value = "10_OLMOS"
$ws.Cells.Item(115,4).Value = 61
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "08:52:26"
$ws.Cells.Item(116,2).Value = "10:10"
$ws.Cells.Item(116,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(116,4).Value = 78
$ws.Cells.Item(116,5).Value = "LP1912"
$ws.Cells.Item(117,1).Value = "08:39:38"
$ws.Cells.Item(117,2).Value = "10:12"
$ws.Cells.Item(117,3).Value = "15_ABASTO"
$ws.Cells.Item(117,4).Value = 93
$ws.Cells.Item(117,5).Value = "LP1912"
$ws.Cells.Item(118,1).Value = "08:52:26"
$ws.Cells.Item(118,2).Value = "10:21"
$ws.Cells.Item(118,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(118,4).Value = 89
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Cells.Item(119,1).Value = "08:52:26"
$ws.Cells.Item(119,2).Value = "10:26"
$ws.Cells.Item(119,3).Value = "215A_EL PATO"
$ws.Cells.Item(119,4).Value = 94
$ws.Cells.Item(119,5).Value = "LP1912"
$ws.Cells.Item(120,1).Value = "08:52:26"
$ws.Cells.Item(120,2).Value = "10:42"
$ws.Cells.Item(120,3).Value = "17_ROMERO"
$ws.Cells.Item(120,4).Value = 110
$ws.Cells.Item(120,5).Value = "LP1912"
$ws.Cells.Item(121,1).Value = "08:52:26"
$ws.Cells.Item(121,2).Value = "10:43"
$ws.Cells.Item(121,3).Value = "14_ABASTO"
$ws.Cells.Item(121,4).Value = 111
$ws.Cells.Item(121,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:52:26"
$ws.Cells.Item(24,1).Value = "08:52:26"
$ws.Cells.Item(24,2).Value = "09:01"
$ws.Cells.Item(24,3).Value = "215A_EL PATO"
$ws.Cells.Item(24,4).Value = 9
$ws.Cells.Item(24,5).Value = "LP1912"
$ws.Cells.Item(26,1).Value = "08:52:26"
$ws.Cells.Item(26,2).Value = "09:42"
$ws.Cells.Item(26,3).Value = "215C_EL PATO"
$ws.Cells.Item(26,4).Value = 50
$ws.Cells.Item(26,5).Value = "LP1912"
$ws.Cells.Item(27,1).Value = "08:52:26"
$ws.Cells.Item(27,2).Value = "10:26"
$ws.Cells.Item(27,3).Value = "215A_EL PATO"
$ws.Cells.Item(27,4).Value = 94
$ws.Cells.Item(27,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:52:26"
$ws.Cells.Item(28,1).Value = "08:52:26"
$ws.Cells.Item(28,2).Value = "09:09"
$ws.Cells.Item(28,3).Value = "215D_LA PLATA"
$ws.Cells.Item(28,4).Value = 17
$ws.Cells.Item(28,5).Value = "L6203"
$ws.Cells.Item(30,1).Value = "08:52:26"
$ws.Cells.Item(30,2).Value = "10:03"
$ws.Cells.Item(30,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(30,4).Value = 71
$ws.Cells.Item(30,5).Value = "L6173"
